{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------------\n// 1) Insert three new sub-bullets right before the \"Can create project\"\n//    bullet. They sit one level deeper (ilvl 4) than \"Can create project\"\n//    (ilvl 3), which is what insertParagraph(..., before) inherits by\n//    default, so the level is bumped explicitly afterwards.\n// ---------------------------------------------------------------------------\nlet createProject = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Can create project\") {\n    createProject = paragraphs.items[i];\n    break;\n  }\n}\n\nconst subBullets = [\n  \"Default start/end dates are blank\",\n  \"Can clear a non-blank date and save it\",\n  \"Blank dates are set to null \",\n];\nfor (const text of subBullets) {\n  const p = createProject.insertParagraph(text, Word.InsertLocation.before);\n  p.listItemOrNullObject.level = 4;\n}\nawait context.sync();\n\n// ---------------------------------------------------------------------------\n// 2) The old \"Can edit project\" bullet becomes \"Can enter blank for dates\n//    and save it\", and a fresh \"Can edit project\" bullet is inserted right\n//    after it (inheriting the same ilvl 3 indent level).\n// ---------------------------------------------------------------------------\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nlet editProject = null;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text.trim() === \"Can edit project\") {\n    editProject = paragraphs2.items[i];\n    break;\n  }\n}\n\neditProject.insertParagraph(\"Can edit project\", Word.InsertLocation.after);\neditProject.insertText(\"Can enter blank for dates and save it\", Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------------\n// 3) Word's \"_GoBack\" bookmark tracked the last edit location (previously on\n//    the \"Can view project details, as user or manager\" bullet). Since the\n//    edit now happens on the \"Can enter blank for dates and save it\"\n//    bullet, move that marker there (remove the old one first, since\n//    bookmark names must stay unique).\n// ---------------------------------------------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nconst endOfEdit = editProject.getRange(\"End\");\nendOfEdit.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# 1) Insert three new sub-bullets right before the \"Can create project\" bullet:\n#    \"Default start/end dates are blank\"\n#    \"Can clear a non-blank date and save it\"\n#    \"Blank dates are set to null \"\n#    These sit one level deeper (ilvl 4 => ListLevelNumber 5) than\n#    \"Can create project\" (ilvl 3 => ListLevelNumber 4).\n# ---------------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq \"Can create project\") {\n        $p.Range.InsertParagraphBefore()\n        $p.Range.Text = \"Default start/end dates are blank\"\n        $p.Range.ListFormat.ListLevelNumber = 5\n\n        $p.Range.InsertParagraphAfter()\n        $p2 = $p.Next()\n        $p2.Range.Text = \"Can clear a non-blank date and save it\"\n        $p2.Range.ListFormat.ListLevelNumber = 5\n\n        $p2.Range.InsertParagraphAfter()\n        $p3 = $p2.Next()\n        $p3.Range.Text = \"Blank dates are set to null \"\n        $p3.Range.ListFormat.ListLevelNumber = 5\n\n        break\n    }\n}\n\n# ---------------------------------------------------------------------------\n# 2) The old \"Can edit project\" bullet becomes \"Can enter blank for dates and\n#    save it\", and a fresh \"Can edit project\" bullet is inserted right after\n#    it (at the same indent level as the old one, ilvl 3).\n# ---------------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq \"Can edit project\") {\n        $p.Range.Text = \"Can enter blank for dates and save it\"\n\n        $p.Range.InsertParagraphAfter()\n        $pNew = $p.Next()\n        $pNew.Range.Text = \"Can edit project\"\n        $pNew.Range.ListFormat.ListLevelNumber = 4\n\n        break\n    }\n}\n\n# ---------------------------------------------------------------------------\n# 3) Word's \"_GoBack\" bookmark tracked the last edit location (previously on\n#    the \"Can view project details, as user or manager\" bullet). Since the\n#    edit now happens on the \"Can enter blank for dates and save it\" bullet,\n#    move that marker there. A zero-length range right at the end of a\n#    paragraph's text needs a real run boundary to anchor to, so a throwaway\n#    placeholder character is inserted, bookmarked against, then removed.\n# ---------------------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq \"Can enter blank for dates and save it\") {\n        $full = $p.Range.Duplicate()\n        $full.MoveEnd(1, -1) | Out-Null\n        $endPos = $full.End\n\n        $full.InsertAfter(\"X\")\n        $anchor = $d.Range($endPos, $endPos)\n        $d.Bookmarks.Add(\"_GoBack\", $anchor)\n\n        $placeholder = $d.Range($endPos, $endPos + 1)\n        $placeholder.Delete()\n        break\n    }\n}\n"}
